# Rename the "CollectionMethods" sheet to "Collection_Methods" to match the
# new naming standard (the _FilterDatabase defined name tracks the sheet
# name automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CollectionMethods")
$ws.Name = "Collection_Methods"
